# Actualización automática 2025-07-21 09:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets("VENTAS POR GRUPO")

# Widen column J (10) from 9 to 11 display units.
# (ColumnWidth 10.17 rounds/snaps internally to a stored width of exactly 11.)
$wsGrupo.Columns.Item(10).ColumnWidth = 10.17

# LED sales value for this advisor/client appears in column J (LED).
$wsGrupo.Range("J18").Value = 73.48

# Update the "x de 30" tally text for column J.
$wsGrupo.Range("J32").Value = "1 de 30"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets("VENTA MENSUAL")

# July (julio) sales updated for this advisor/client row.
$wsMensual.Range("F18").Value = 1788.22

# July total updated accordingly.
$wsMensual.Range("F32").Value = 4611.43

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets("CUMPLIMIENTO MENSUAL")

# LED row (row 8): VENTA, POR CUMPLIR and CUMPLIMIENTO updated.
$wsCumplimiento.Range("D8").Value = 73.48
$wsCumplimiento.Range("E8").Value = 226.52
$wsCumplimiento.Range("F8").Value = 0.2449333333333333

# TOTAL row (row 18): VENTA, POR CUMPLIR and CUMPLIMIENTO updated.
$wsCumplimiento.Range("D18").Value = 4601.150000000001
$wsCumplimiento.Range("E18").Value = 29333.56607548726
$wsCumplimiento.Range("F18").Value = 0.1355882863367653
